# RBA v2.5 - Atualizacao da Tela
# Replace the placeholder "Tre"-family tokens with "Qwer"-family tokens
# throughout the convocation document: the bold salutation run in the
# body, and the letterhead block in the page header (region name,
# address line, CEP/telephone/e-mail lines).

$d = $word.ActiveDocument

# wdReplaceOne / wdReplaceAll constants used with Find.Execute's Replace arg
$wdReplaceOne = 1
$wdReplaceAll = 2

# --- 1. Document body: the bold "TERE" in "A TERE, vem por meio desta ..." ---
$d.Content.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", $wdReplaceAll)

# --- 2. Page header: letterhead block ---
# wdPrimaryHeaderStory = 7
$hdr = $d.StoryRanges(7)

# "DIRETORIA DE ENSINO REGIAO TRE" -> "... QWER"
$hdr.Find.Execute("TRE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", $wdReplaceOne)

# own paragraph "TERE - DEP." -> "QWER - DEP."
$hdr.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", $wdReplaceOne)

# Address line: "Tre, no Tre - Tre - Tre - Tre" -> "Qwer, no Qwer - Qewr - Qewr - Qwer"
$hdr.Find.Execute("Tre", $true, $false, $false, $false, $false, $true, 1, $false, "Qwer", $wdReplaceOne)
$hdr.Find.Execute("Tre", $true, $false, $false, $false, $false, $true, 1, $false, "Qwer", $wdReplaceOne)
$hdr.Find.Execute("Tre", $true, $false, $false, $false, $false, $true, 1, $false, "Qewr", $wdReplaceOne)
$hdr.Find.Execute("Tre", $true, $false, $false, $false, $false, $true, 1, $false, "Qewr", $wdReplaceOne)
$hdr.Find.Execute("Tre", $true, $false, $false, $false, $false, $true, 1, $false, "Qwer", $wdReplaceOne)

# CEP / Tel / Email lines (lower-case "tre")
$hdr.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, "qwer", $wdReplaceOne)
$hdr.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, "qwer", $wdReplaceOne)
$hdr.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, "qwer", $wdReplaceOne)

Write-Host "Header now reads:" $d.StoryRanges(7).Text
